$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 (B39): shorten/replace the Pengu swim intro text
$ws.Range("B39").Value = "Help Pengu swim towards the land!"

# New row 46: level_7_intro_1 dialog
$ws.Range("A46").Value = "level_7_intro_1"
$ws.Range("B46").Value = "A boulder is blocking Pengu's path!"
$ws.Range("B46").VerticalAlignment = -4108

# New row 47: level_7_intro_2 dialog
$ws.Range("A47").Value = "level_7_intro_2"
$ws.Range("B47").Value = "It is time to unleash the most potent of Pengu's fractional powers!"
$ws.Range("B47").VerticalAlignment = -4108

# Update the view to reflect the new extent of data (selection on the last edited cell)
$ws.Range("B47").Select() | Out-Null
